$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date value for rows 2-6 from 45224 to 45233
$ws.Range("C2:C6").Value = 45233
